$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.471.24'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.727.20'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9971'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.88'
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9977'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4894'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2607'
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06206'
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').Value = '1.728.00'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06989'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.62'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.531'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6016'
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.20'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9974'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '26.452.25'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9972'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').Value = '1.943.87'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.467'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.513'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.109'
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.45'
$ws.Range('E25').Value = '  -2.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.32'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.749'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.71'
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.918'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08013'
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.644'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04494'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9966'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.598'
$ws.Range('E35').Value = '  -0.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.003'
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6231'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9353'
$ws.Range('E38').Value = '  +3.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.993'
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9972'
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01481'
$ws.Range('E42').Value = '  -1.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.79'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.411'
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3853'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.908'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05369'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.53'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.722'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.43'
$ws.Range('E51').Value = '  -0.52%  '
